$wb = $excel.ActiveWorkbook

# Worksheets: 1 = AppControl, 2 = smoke, 3 = regression, 4 = miniregression, 5 = DeviceName
$wsAppControl = $wb.Worksheets.Item(1)
$wsSmoke = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# smoke sheet: refresh the COBT / sensor test case names and add a new one
# ---------------------------------------------------------------------------
$wsSmoke.Range("A17").Value = "ZestIOT_2268_Verifying_COBT_For_DIALCelebi_User"
$wsSmoke.Range("A18").Value = "ZestIOT_2268_Verifying_COBT_For_GMR_HYD_AISATS_User"
$wsSmoke.Range("A19").Value = "ZestIOT_2268_Verifying_COBT_For_GMR_HYD_SG_User"
$wsSmoke.Range("A20").Value = "ZestIOT_2293_GMR_HYD_Sensor_And_Scheduled_data_Validation"

# Copy B20 (carries the run-flag style) down into the newly added row 21
$wsSmoke.Range("B20").Copy($wsSmoke.Range("B21")) | Out-Null
$wsSmoke.Range("A21").Value = "ZestIOT_2294_GMR_HYD_SensorATD_And_Scheduled_data_Validation"
$wsSmoke.Range("B21").Value = "Y"

# ---------------------------------------------------------------------------
# AppControl sheet: add the email address + hyperlink in B25 (Email ID row)
# ---------------------------------------------------------------------------
$wsAppControl.Range("B25").Value = "stiyyagura@enhops.com"

$origStyle = $wsAppControl.Range("B25").Style
$wsAppControl.Hyperlinks.Add($wsAppControl.Range("B25"), "mailto:stiyyagura@enhops.com") | Out-Null
$wsAppControl.Range("B25").Style = $origStyle

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping to match the saved view state
# ---------------------------------------------------------------------------
$wsAppControl.Range("A26").Select() | Out-Null

$wsSmoke.Activate() | Out-Null
$wsSmoke.Range("A18").Select() | Out-Null
